# Updates cryptos list prices / 1h volume percentages (and restores the
# correct coin/link ordering for two swapped row pairs), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: several Price values (column D) look like plain decimal numbers
# (e.g. "253.90"); a leading apostrophe is used so Excel stores them as
# text (matching the workbook's existing text-based Price column) rather
# than silently coercing them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.184.36"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "3.389.83"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'253.90"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "'662.64"
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  -2.65%  "
$ws.Range("D9").Value = "'1.04"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "3.388.53"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "'41.61"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").Value = "97.825.32"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("D16").Value = "'0.0000256"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("D17").Value = "4.015.04"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("E18").Value = "  -2.97%  "
$ws.Range("D19").Value = "3.384.34"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").Value = "'18.02"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'0.528"
$ws.Range("E21").Value = "  -8.11%  "
$ws.Range("D22").Value = "'10.96"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'3.45"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'512.01"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "'6.99"
$ws.Range("E25").Value = "  +7.56%  "
$ws.Range("D26").Value = "'0.0000202"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").Value = "'96.67"
$ws.Range("E27").Value = "  -4.29%  "
$ws.Range("D28").Value = "'12.36"
$ws.Range("E28").Value = "  -4.94%  "
$ws.Range("D29").Value = "3.568.54"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("E30").Value = "  -2.75%  "
$ws.Range("D31").Value = "'0.143"
$ws.Range("E31").Value = "  -4.79%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  -6.17%  "
$ws.Range("D34").Value = "'2.60"
$ws.Range("E34").Value = "  +8.11%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "'0.560"
$ws.Range("E36").Value = "  -3.70%  "
$ws.Range("D37").Value = "'28.95"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("D38").Value = "'8.00"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "'1.51"
$ws.Range("E39").Value = "  -1.81%  "
$ws.Range("D40").Value = "'535.21"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'24.42"
$ws.Range("E43").Value = "  -1.13%  "
$ws.Range("D44").Value = "'0.857"
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "'0.0428"
$ws.Range("E45").Value = "  -1.04%  "
$ws.Range("D46").Value = "'1.72"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").Value = "'3.68"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.26"
$ws.Range("E48").Value = "  +6.19%  "
$ws.Range("D49").Value = "'5.62"
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("D50").Value = "'56.09"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").Value = "'8.59"
$ws.Range("E51").Value = "  -6.03%  "
